# rabbit-test.xlsx: "now parses host string"
#
# The RTD server used to be handed a bare hostname ("localhost"); it now
# parses a full connection-string-like host ("amqp://guest:guest@localhost"),
# and the RTD() calls gained an extra (blank) topic parameter right after
# EXCHANGE to carry it. Also bump the "added by code" timestamp cell and
# move the active selection down one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- RABBIT_URI (named range -> Sheet2!$D$2): "localhost" -> full URI ---
# Duplicate D2 (value + formatting) into H2 first so H2 picks up the same
# style (s="5") the old E2 helper cell had, then overwrite both with the
# new URI text and drop the now-redundant E2 cell entirely.
$ws.Range("D2").Copy($ws.Range("H2")) | Out-Null
$ws.Range("D2").Value = "amqp://guest:guest@localhost"
$ws.Range("H2").Value = "amqp://guest:guest@localhost"
$ws.Range("E2").Clear() | Out-Null

# --- "This text was added by using code ..." timestamp banner ---
$ws.Range("A1").Value = "This text was added by using code 6/14/2018 1:41:10 PM"

# --- Table1 calculated columns: RTD(...) gains an extra blank arg ---
# (progId,,RABBIT_URI,EXCHANGE,Table1[...]  ->  progId,,RABBIT_URI,EXCHANGE,,Table1[...])
$fieldFormula = '=RTD(progId,,RABBIT_URI,EXCHANGE,,Table1[[#This Row],[routingKey]],Table1[[#Headers],[FIELD]])'
$blankFormula = '=RTD(progId,,RABBIT_URI,EXCHANGE,,Table1[[#This Row],[routingKey]],IF(Table1[[#Headers],[Blank]]="Blank","",Table1[[#Headers],[Blank]]) )'

for ($row = 5; $row -le 12; $row++) {
    $ws.Range("C$row").Formula = $fieldFormula
    $ws.Range("D$row").Formula = $blankFormula
}

# --- selection moved from D6 to D7 ---
$ws.Range("D7").Select() | Out-Null
